$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values (rows 2-10, columns B-O), derived from target diff.
$data = @{
    2  = @(38,2,23,3,25,4,14,5,6,6,22,6,20,6)
    3  = @(25,2,24,3,18,4,6,5,7,6,14,6,19,6)
    4  = @(23,2,17,3,15,4,15,5,46,5,13,6,39,6)
    5  = @(25,2,18,3,17,4,13,5,4,6,16,6,18,6)
    6  = @(27,2,9,3,16,4,3,5,4,6,15,6,25,6)
    7  = @(33,2,6,3,8,3,2,4,14,4,16,6,23,6)
    8  = @(30,2,46,2,15,2,6,2,8,2,12,3,33,6)
    9  = @(40,2,35,2,10,2,6,1,15,1,5,2,11,6)
    10 = @(25,2,31,2,27,1,5,1,4,0,1,1,1,4)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2   # column B = 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}

# Update the selected cell in the sheet view
$ws.Range("O8").Select()
